$wb = $excel.ActiveWorkbook

# --- SF1_unrooted: add new row 14 (AM_GG) ---
$wsSF1u = $wb.Worksheets.Item("SF1_unrooted")
$wsSF1u.Range("B13").Copy()
$wsSF1u.Range("B14").PasteSpecial(-4122)
$wsSF1u.Range("A14").Value = "AM_GG"
$wsSF1u.Range("B14").Value = -4871.0361810000004
$wsSF1u.Range("C14").Formula = "=2*(B14-B2)"
$wsSF1u.Range("C14").Select() | Out-Null

# --- WT1_unrooted: add new row 14 (AM_GG) ---
$wsWT1u = $wb.Worksheets.Item("WT1_unrooted")
$wsWT1u.Range("B13").Copy()
$wsWT1u.Range("B14").PasteSpecial(-4122)
$wsWT1u.Range("A14").Value = "AM_GG"
$wsWT1u.Range("B14").Value = -3697.1109809999998
$wsWT1u.Range("C14").Formula = "=2*(B14-B2)"
$wsWT1u.Range("C14").Select() | Out-Null
